$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 505
$ws.Range("I2").Value = 1246
$ws.Range("J2").Value = 5517
$ws.Range("L2").Value = 1558
$ws.Range("M2").Value = 95
$ws.Range("N2").Value = 942
$ws.Range("O2").Value = 3
$ws.Range("Q2").Value = 7
$ws.Range("R2").Value = 80
$ws.Range("S2").Value = 609
$ws.Range("T2").Value = 995
$ws.Range("U2").Value = 68
$ws.Range("V2").Value = 8740
$ws.Range("X2").Value = 8760
$ws.Range("Y2").Value = 12
$ws.Range("Z2").Value = 132
$ws.Range("AA2").Value = 51
